$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "2023-10-07 23:12:07 36_12_6534199"
$ws.Range("B6").Value = "What do you call a cow with no legs?"
$ws.Range("C6").Value = "Ground beef!"
